# -----------------------------------------------------------------------
# Se implemento backpropagation con inercia: se agregan las hojas
# "Red con inercia" (pruebas variando mu) y "Red con Cross entropy"
# (hoja nueva, vacia, para el siguiente experimento).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Crear las dos hojas nuevas, en orden, justo despues de "Red basica"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Red con inercia"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Red con Cross entropy"

# ---------------------------------------------------------------------
# 2) "Red con inercia": copiar la tabla base (A1:F7, encabezados +
#    filas) de "Red basica" para heredar exactamente el mismo formato,
#    luego copiar la columna de "Descripcion" (antes G) a H, y volver a
#    copiar esa misma columna a G para heredar el estilo antes de
#    sobrescribir sus valores con la nueva columna "mu".
# ---------------------------------------------------------------------
[void]$ws1.Range("A1:F7").Copy($ws2.Range("A1:F7"))
[void]$ws1.Range("G1:G7").Copy($ws2.Range("H1:H7"))
[void]$ws1.Range("G1:G7").Copy($ws2.Range("G1:G7"))

# eta se mantuvo fijo en 7.5 (el mejor valor obtenido con la red basica)
# en todas las pruebas de esta hoja.
$ws2.Range("F2:F7").Value = 7.5

# Encabezado de la nueva columna G ("mu")
$ws2.Range("G1").Value = "mu"
$ws2.Range("H1").Value = "Descripción"

# Valores de mu ensayados y su descripcion de resultado
$ws2.Range("G2").Value = 0.4
$ws2.Range("H2").Value = "Precisión en datos de validación del 73.06%."

$ws2.Range("G3").Value = 0.5
$ws2.Range("H3").Value = "Precisión en datos de validación del 92.68%"

$ws2.Range("G4").Value = 0.6
$ws2.Range("H4").Value = "Precisión en datos de validación del 92.95%."

$ws2.Range("G5").Value = 0.7
$ws2.Range("H5").Value = "Precisión en datos de validación del 83.41%."

$ws2.Range("G6").Value = 0.8
$ws2.Range("H6").Value = "Precisión en datos de validación del 92.91%."

$ws2.Range("G7").Value = 0.9
$ws2.Range("H7").Value = "Precisión en datos de validación del 93.40%."

# Anchos de columna G y H
$ws2.Columns.Item(7).ColumnWidth = 8.3
$ws2.Columns.Item(8).ColumnWidth = 85.45

# ---------------------------------------------------------------------
# 3) "Red basica": se quitan las filas 8:10 (ya no se usan) y se mueve
#    la seleccion a G9.
# ---------------------------------------------------------------------
[void]$ws1.Rows("8:10").Delete()
[void]$ws1.Range("G9").Select()

# ---------------------------------------------------------------------
# 4) Seleccion final / hoja activa: "Red con inercia"
# ---------------------------------------------------------------------
[void]$ws2.Range("H8").Select()
$ws2.Activate()
